# Apply "carls newest scraper" update:
#   - Drop the "Norm, Typ" column (A) and several unused columns
#     (Ritningsnummer, Position, Beteckning, Kompletterande Information
#     ovrigt, Ref annan, Historiskt Varumarke, Historiskt inkopsreferens,
#     Forpackning) so only Varumarke..Typbeteckning, Enhet, SSG-notering,
#     E-nummer and RSK-nummer remain (A:I).
#   - Drop the (now stale) AutoFilter / sort state.
#   - Resize the new "SSG-notering" column (G) and refresh the selection.
#   - Shrink the _FilterDatabase defined name to the new used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the unwanted columns. Work right-to-left so earlier deletes don't
# shift the column letters used by later ones.
$ws.Columns("Q").Delete()        # blank spacer column before E-nummer/RSK-nummer
$ws.Columns("O").Delete()        # Forpackning
$ws.Range("G1:M1").EntireColumn.Delete()   # Ritningsnummer..Historiskt inkopsreferens
$ws.Columns("A").Delete()        # Norm, Typ

# The old AutoFilter (and its cached sort state) referenced the old, much
# wider/deeper range - just drop it.
$ws.AutoFilterMode = $false

# Shrink the hidden _xlnm._FilterDatabase name down to the surviving columns.
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Sheet1!_FilterDatabase") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$I`$21217"
  }
}

# Widen the "SSG-notering" column (now column G) to fit its content.
$ws.Columns("G").ColumnWidth = 28.6

# Refresh the selection/scroll position to the new RSK-nummer-ish column.
$ws.Range("H1:H1048576").Select()
